$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.088.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.963.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4973"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4220"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09115"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.096"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.60%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.972.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.894"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.48%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.429"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06682"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.17%  "
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.925"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.107.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.291"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.233"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.268"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.040"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09844"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.532"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.812"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.698"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02428"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.028"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06348"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.285"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6447"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.181"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.471"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06865"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.109"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.17%  "
